$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 17998.166
$ws.Range("I70").Value = 1266.6666
$ws.Range("J70").Value = 34729.668
$ws.Range("K70").Value = 3799.9998
$ws.Range("L70").Value = 104189.004
$ws.Range("M70").Value = -3529.9998
$ws.Range("N70").Value = -104729.004
$ws.Range("H73").Value = 17998.166
$ws.Range("I73").Value = 1266.6666
$ws.Range("J73").Value = 34729.668
$ws.Range("K73").Value = 3799.9998
$ws.Range("L73").Value = 104189.004
$ws.Range("M73").Value = -2863.9998
$ws.Range("N73").Value = -106061.004
$ws.Range("H80").Value = 877.4231
$ws.Range("I80").Value = 815.8125
$ws.Range("K80").Value = 2447.4375
$ws.Range("M80").Value = -1449.4375
$ws.Range("H83").Value = 877.4231
$ws.Range("I83").Value = 815.8125
$ws.Range("K83").Value = 7342.3125
$ws.Range("M83").Value = -2350.3125
$ws.Range("H86").Value = 2498.5
$ws.Range("I86").Value = 998.2
$ws.Range("K86").Value = 998.2
$ws.Range("M86").Value = 124.8
$ws.Range("H89").Value = 2498.5
$ws.Range("I89").Value = 998.2
$ws.Range("K89").Value = 4991
$ws.Range("M89").Value = 625
$ws.Range("H107").Value = 992.2
$ws.Range("I107").Value = 347.33334
$ws.Range("K107").Value = 347.33334
$ws.Range("M107").Value = 1572.66666
$ws.Range("H112").Value = 3999.5
$ws.Range("I112").Value = 1100
$ws.Range("J112").Value = 4192.8
$ws.Range("K112").Value = 3300
$ws.Range("L112").Value = 12578.4
$ws.Range("M112").Value = -2192
$ws.Range("N112").Value = -14794.4
$ws.Range("H132").Value = 1505.8948
$ws.Range("I132").Value = 1200.0714
$ws.Range("J132").Value = 2362.2
$ws.Range("K132").Value = 3600.2142
$ws.Range("L132").Value = 7086.599999999999
$ws.Range("M132").Value = -1070.2142
$ws.Range("N132").Value = -12146.6
$ws.Range("H135").Value = 715.4286
$ws.Range("I135").Value = 668
$ws.Range("K135").Value = 6012
$ws.Range("M135").Value = -3477
$ws.Range("H137").Value = 1530.625
$ws.Range("I137").Value = 1536.4286
$ws.Range("J137").Value = 1490
$ws.Range("K137").Value = 4609.2858
$ws.Range("L137").Value = 4470
$ws.Range("M137").Value = -2059.2858
$ws.Range("N137").Value = -9570
$ws.Range("H141").Value = 3381.8845
$ws.Range("I141").Value = 2695.5789
$ws.Range("J141").Value = 5244.7144
$ws.Range("K141").Value = 8086.736699999999
$ws.Range("L141").Value = 15734.1432
$ws.Range("M141").Value = -2906.736699999999
$ws.Range("N141").Value = -26094.1432

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5435.029
$ws.Range("I32").Value = 3740.5652
$ws.Range("J32").Value = 8682.75
$ws.Range("K32").Value = 3740.5652
$ws.Range("L32").Value = 8682.75
$ws.Range("M32").Value = -3453.5652
$ws.Range("N32").Value = -9256.75
$ws.Range("H61").Value = 4783.8
$ws.Range("I61").Value = 2653
$ws.Range("K61").Value = 2653
$ws.Range("M61").Value = -2441
$ws.Range("H74").Value = 1283.8695
$ws.Range("I74").Value = 1110.5625
$ws.Range("K74").Value = 1110.5625
$ws.Range("M74").Value = -236.5625
$ws.Range("H77").Value = 1283.8695
$ws.Range("I77").Value = 1110.5625
$ws.Range("K77").Value = 5552.8125
$ws.Range("M77").Value = -1184.8125
$ws.Range("H122").Value = 633
$ws.Range("I122").Value = 633
$ws.Range("K122").Value = 1899
$ws.Range("M122").Value = 551
$ws.Range("H136").Value = 4783.8
$ws.Range("I136").Value = 2653
$ws.Range("K136").Value = 7959
$ws.Range("M136").Value = -5409

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1336
$ws.Range("I16").Value = 1037.8
$ws.Range("K16").Value = 1037.8
$ws.Range("M16").Value = -750.8
$ws.Range("H31").Value = 1834.9697
$ws.Range("I31").Value = 1379.238
$ws.Range("K31").Value = 1379.238
$ws.Range("M31").Value = -1084.238
$ws.Range("H34").Value = 1834.9697
$ws.Range("I34").Value = 1379.238
$ws.Range("K34").Value = 1379.238
$ws.Range("M34").Value = -1177.238
$ws.Range("H86").Value = 250002320
$ws.Range("I86").Value = 1000000000
$ws.Range("J86").Value = 3099.3333
$ws.Range("K86").Value = 1000000000
$ws.Range("L86").Value = 3099.3333
$ws.Range("M86").Value = -999998877
$ws.Range("N86").Value = -5345.3333
$ws.Range("H89").Value = 250002320
$ws.Range("I89").Value = 1000000000
$ws.Range("J89").Value = 3099.3333
$ws.Range("K89").Value = 5000000000
$ws.Range("L89").Value = 15496.6665
$ws.Range("M89").Value = -4999994384
$ws.Range("N89").Value = -26728.6665
$ws.Range("H94").Value = 1007.8
$ws.Range("J94").Value = 1020.875
$ws.Range("L94").Value = 1020.875
$ws.Range("N94").Value = -1922.875
$ws.Range("H107").Value = 479.7647
$ws.Range("I107").Value = 399.69232
$ws.Range("J107").Value = 740
$ws.Range("K107").Value = 399.69232
$ws.Range("L107").Value = 740
$ws.Range("M107").Value = 1520.30768
$ws.Range("N107").Value = -4580
$ws.Range("H113").Value = 1336
$ws.Range("I113").Value = 1037.8
$ws.Range("K113").Value = 1037.8
$ws.Range("M113").Value = 1132.2
$ws.Range("H132").Value = 1854.8864
$ws.Range("I132").Value = 1266.2258
$ws.Range("K132").Value = 3798.6774
$ws.Range("M132").Value = -1268.6774
$ws.Range("H134").Value = 1698.9706
$ws.Range("I134").Value = 1481.1613
$ws.Range("K134").Value = 4443.4839
$ws.Range("M134").Value = -1908.4839

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 301
$ws.Range("J92").Value = 301.5
$ws.Range("L92").Value = 904.5
$ws.Range("N92").Value = -3400.5
$ws.Range("H97").Value = 980
$ws.Range("I97").Value = 0
$ws.Range("J97").Value = 980
$ws.Range("K97").Value = 0
$ws.Range("L97").Value = 2940
$ws.Range("M97").ClearContents()
$ws.Range("N97").Value = -3932
$ws.Range("H131").Value = 770.3099999999999
$ws.Range("I131").Value = 493.7143
$ws.Range("J131").Value = 791.129
$ws.Range("K131").Value = 1481.1429
$ws.Range("L131").Value = 2373.387
$ws.Range("M131").Value = 3558.8571
$ws.Range("N131").Value = -12453.387

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 65.882355
$ws.Range("I2").Value = 13.8
$ws.Range("J2").Value = 87.583336
$ws.Range("K2").Value = 13.8
$ws.Range("L2").Value = 87.583336
$ws.Range("M2").Value = 99.2
$ws.Range("N2").Value = -313.583336
$ws.Range("H80").Value = 1600
$ws.Range("I80").Value = 300
$ws.Range("K80").Value = 300
$ws.Range("M80").Value = 698
$ws.Range("H83").Value = 1600
$ws.Range("I83").Value = 300
$ws.Range("K83").Value = 1500
$ws.Range("M83").Value = 3492

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H4").Value = 7665.0625
$ws.Range("I4").Value = 6211.75
$ws.Range("J4").Value = 12025
$ws.Range("K4").Value = 6211.75
$ws.Range("L4").Value = 12025
$ws.Range("M4").Value = -6098.75
$ws.Range("N4").Value = -12251
$ws.Range("H28").Value = 7665.0625
$ws.Range("I28").Value = 6211.75
$ws.Range("J28").Value = 12025
$ws.Range("K28").Value = 6211.75
$ws.Range("L28").Value = 12025
$ws.Range("M28").Value = -5979.75
$ws.Range("N28").Value = -12489
$ws.Range("H37").Value = 7665.0625
$ws.Range("I37").Value = 6211.75
$ws.Range("J37").Value = 12025
$ws.Range("K37").Value = 6211.75
$ws.Range("L37").Value = 12025
$ws.Range("M37").Value = -6104.75
$ws.Range("N37").Value = -12239
$ws.Range("H132").Value = 3117.56
$ws.Range("I132").Value = 2495.8572
$ws.Range("J132").Value = 3908.818
$ws.Range("K132").Value = 7487.571599999999
$ws.Range("L132").Value = 11726.454
$ws.Range("M132").Value = -4957.571599999999
$ws.Range("N132").Value = -16786.454
$ws.Range("H136").Value = 2841.7273
$ws.Range("I136").Value = 2049
$ws.Range("K136").Value = 6147
$ws.Range("M136").Value = -3597

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H26").Value = 25000000
$ws.Range("J26").Value = 25000000
$ws.Range("L26").Value = 25000000
$ws.Range("N26").Value = -25000586
$ws.Range("H70").Value = 29159.4
$ws.Range("J70").Value = 29159.4
$ws.Range("L70").Value = 29159.4
$ws.Range("N70").Value = -29789.4
$ws.Range("H73").Value = 29159.4
$ws.Range("J73").Value = 29159.4
$ws.Range("L73").Value = 29159.4
$ws.Range("N73").Value = -31343.4
$ws.Range("H122").Value = 87726.22
$ws.Range("I122").Value = 130846.164
$ws.Range("J122").Value = 1486.3334
$ws.Range("K122").Value = 392538.492
$ws.Range("L122").Value = 4459.0002
$ws.Range("M122").Value = -390088.492
$ws.Range("N122").Value = -9359.0002
$ws.Range("H126").Value = 1375.8667
$ws.Range("I126").Value = 1516.9166
$ws.Range("J126").Value = 811.6667
$ws.Range("K126").Value = 4550.7498
$ws.Range("L126").Value = 2435.0001
$ws.Range("M126").Value = -2080.7498
$ws.Range("N126").Value = -7375.0001
$ws.Range("H132").Value = 1970.8334
$ws.Range("I132").Value = 1334.2727
$ws.Range("K132").Value = 4002.8181
$ws.Range("M132").Value = -1472.8181
